# Edit: add Galar-region Pokemon rows (Scorbunny line through Drednaw)
# to the "Tabela Pokemon" worksheet, mirroring the author's commit
# "Projeto de Regressao Logistica, Advertising" data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 813 (Rillaboom) was missing its "Lendario" flag; fill it in to match
# every other row in the table.
$ws.Range("R813").Value = "Nao"

# Each row: Nome, Numero, Tipo 1, Tipo 2, Altura m, Peso Kg, Catch Rate,
#           Male, Female, Hp, Atk, Def, SpA, SpD, Spe, Lendario
$newRows = @(
    @("Scorbunny", 813, "Fire", $null, 0.3, 4.5, 45, 87.5, 12.5, 50, 71, 40, 40, 40, 69, "Nao"),
    @("Raboot", 814, "Fire", $null, 0.6, 9, 45, 87.5, 12.5, 65, 86, 60, 55, 60, 94, "Nao"),
    @("Cinderace", 815, "Fire", $null, 1.4, 33, 45, 87.5, 12.5, 80, 116, 75, 65, 75, 119, "Nao"),
    @("Sobble", 816, "Water", $null, 0.3, 4, 45, 87.5, 12.5, 50, 40, 40, 70, 40, 70, "Nao"),
    @("Drizzile", 817, "Water", $null, 0.7, 11.5, 45, 87.5, 12.5, 65, 60, 55, 95, 55, 90, "Nao"),
    @("Inteleon", 818, "Water", $null, 1.9, 45.2, 45, 87.5, 12.5, 70, 85, 65, 125, 65, 120, "Nao"),
    @("Skwovet", 819, "Normal", $null, 0.3, 2.5, 255, 50, 50, 70, 55, 55, 35, 35, 25, "Nao"),
    @("Greedent", 820, "Normal", $null, 0.6, 6, 90, 50, 50, 120, 95, 95, 55, 75, 20, "Nao"),
    @("Rookidee", 821, "Flying", $null, 0.2, 1.8, 255, 50, 50, 38, 47, 35, 33, 35, 57, "Nao"),
    @("Corvisquire", 822, "Flying", $null, 0.8, 16, 120, 50, 50, 68, 67, 55, 43, 55, 77, "Nao"),
    @("Corviknight", 823, "Flying", "Steel", 2.2000000000000002, 75, 45, 50, 50, 98, 87, 105, 53, 85, 67, "Nao"),
    @("Blipbug", 824, "Bug", $null, 0.4, 8, 255, 50, 50, 25, 20, 20, 25, 45, 45, "Nao"),
    @("Dottler", 825, "Bug", "Psychic", 0.4, 19.5, 120, 50, 50, 50, 35, 80, 50, 90, 30, "Nao"),
    @("Orbeetle", 826, "Bug", "Psychic", 0.4, 40.799999999999997, 45, 50, 50, 60, 45, 110, 80, 120, 90, "Nao"),
    @("Nickit", 827, "Dark", $null, 0.6, 8.9, 255, 50, 50, 40, 28, 28, 47, 52, 50, "Nao"),
    @("Thievul", 828, "Dark", $null, 1.2, 19.899999999999999, 127, 50, 50, 70, 58, 58, 87, 92, 90, "Nao"),
    @("Gossifleur", 829, "Grass", $null, 0.4, 2.2000000000000002, 190, 50, 50, 40, 40, 60, 40, 60, 10, "Nao"),
    @("Eldegoss", 830, "Grass", $null, 0.5, 2.5, 75, 50, 50, 60, 50, 90, 80, 120, 60, "Nao"),
    @("Wooloo", 831, "Normal", $null, 0.6, 6, 255, 50, 50, 42, 40, 55, 40, 45, 48, "Nao"),
    @("Dubwool", 832, "Normal", $null, 1.3, 43, 127, 50, 50, 72, 80, 100, 60, 90, 88, "Nao"),
    @("Chewtle", 833, "Water", $null, 0.3, 8.5, 255, 50, 50, 50, 64, 50, 38, 38, 44, "Nao"),
    @("Drednaw", 834, "Water", "Rock", 1, 115.5, 75, 50, 50, 90, 115, 90, 48, 68, 74, "Nao")
)

$startRow = 814
$colLetters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","R")

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $colLetters.Count; $c++) {
        $val = $rowData[$c]
        if ($null -ne $val) {
            $ws.Range($colLetters[$c] + $r).Value = $val
        }
    }
    $ws.Range("P" + $r).Formula = "=SUM(J" + $r + ":O" + $r + ")"
    $ws.Range("Q" + $r).Formula = "=AVERAGE(J" + $r + ":O" + $r + ")"
}

$lastRow = $startRow + $newRows.Count - 1

# Move the selection to the first empty row below the new data, matching
# the post-edit cursor position recorded in the workbook.
[void]$ws.Range("A" + ($lastRow + 1)).Select()
